# Database.xlsx edit script
# - Adds a new "SRS830TimeConstants" worksheet with the lock-in amplifier
#   time-constant lookup table.
# - Reorders the worksheet tabs to Magnification, Materials, SRS830TimeConstants.
# - Makes the new sheet the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new sheet at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "SRS830TimeConstants"

# ---------------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "index"
$ws.Range("B1").Value = "timeConstant"
$ws.Range("A1:B1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. Data rows: SRS830 lock-in amplifier time constants, index 0-19.
# ---------------------------------------------------------------------------
$values = @(
    @(0,  0.000001),
    @(1,  0.000003),
    @(2,  0.0001),
    @(3,  0.0003),
    @(4,  0.001),
    @(5,  0.003),
    @(6,  0.01),
    @(7,  0.03),
    @(8,  0.1),
    @(9,  0.3),
    @(10, 1),
    @(11, 3),
    @(12, 10),
    @(13, 30),
    @(14, 100),
    @(15, 300),
    @(16, 1000),
    @(17, 3000),
    @(18, 10000),
    @(19, 30000)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$ws.Range("B2:B21").NumberFormat = "##0.0E+0"

# ---------------------------------------------------------------------------
# 4. Column widths / view.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 6
$ws.Columns.Item(2).ColumnWidth = 13.140625

$ws.Range("D43").Select()

# ---------------------------------------------------------------------------
# 5. Reorder the tabs: Magnification, Materials, SRS830TimeConstants.
# ---------------------------------------------------------------------------
$magnification = $wb.Worksheets.Item("Magnification")
$magnification.Move($wb.Worksheets.Item(1))

# Re-acquire the sheet reference post-move before touching it further.
$magnification = $wb.Worksheets.Item("Magnification")

# Magnification sheet no longer needs to be the selected tab; give it a
# frozen header row (row 1) with C2 selected, matching the refreshed layout,
# and a portrait page setup (mirrors the printer setup Materials already had).
$magnification.Activate()
$magnification.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$magnification.Range("C2").Select()
$magnification.Range("A1:B1").Font.Bold = $true
$magnification.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 6. Make the new sheet the active tab.
# ---------------------------------------------------------------------------
$ws.Activate()
